$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook originally has 4 sheets: 总计, 2022-Q3, 2021-Q2, 2021-Q1
# We need to insert a new "2022-Q4" sheet right after "总计", so the final
# order becomes: 总计, 2022-Q4, 2022-Q3, 2021-Q2, 2021-Q1
# ---------------------------------------------------------------------------

$wsTotal = $wb.Worksheets.Item(1)
$wsQ3    = $wb.Worksheets.Item("2022-Q3")

# Create the new "2022-Q4" sheet as a copy of "2022-Q3" (so it inherits the
# same column layout / cell styles), positioned right after "总计".
$wsQ3.Copy($null, $wsTotal)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# Fill in the "2022-Q4" sheet with the fund-position data (rows 2-9).
# Row 1 (the header) and the base formatting already came from the copy of
# "2022-Q3", so we only need to overwrite the data rows.
# ---------------------------------------------------------------------------

# Make sure rows 4-9 exist with the same per-row styling (border etc.) as the
# already-present rows 2/3 (copied from "2022-Q3") by copying the format of
# row 3 down across rows 4-9.
$wsQ4.Range("A3:H3").Copy()
$wsQ4.Range("A4:H9").PasteSpecial(-4122)

# Columns B (fund code) and D-G (scale/position/ratio/market-value, which are
# stored as plain text strings in the source data) must stay text so that
# leading zeros / exact decimal text are preserved instead of being coerced
# to numbers.
$wsQ4.Range("B2:B9").NumberFormat = "@"
$wsQ4.Range("D2:G9").NumberFormat = "@"
# ...except G9, which is numeric (0) in the source data.
$wsQ4.Range("G9").NumberFormat = "General"

function Set-Q4Row($row, $idx, $code, $name, $scale, $pos, $ratio, $mv, $rank) {
    $wsQ4.Cells.Item($row, 1).Value = $idx
    $wsQ4.Cells.Item($row, 2).Value = $code
    $wsQ4.Cells.Item($row, 3).Value = $name
    $wsQ4.Cells.Item($row, 4).Value = $scale
    $wsQ4.Cells.Item($row, 5).Value = $pos
    $wsQ4.Cells.Item($row, 6).Value = $ratio
    $wsQ4.Cells.Item($row, 7).Value = $mv
    $wsQ4.Cells.Item($row, 8).Value = $rank
}

Set-Q4Row 2 0 "010326" "博时消费创新混合A"     "14.23" "80.27" "2.91" "0.4141" 10
Set-Q4Row 3 1 "013836" "博时时代消费混合A"     "4.75"  "92.95" "6.26" "0.2974" 1
Set-Q4Row 4 2 "004505" "博时新兴消费主题混合A" "4.66"  "87.10" "5.81" "0.2707" 5
Set-Q4Row 5 3 "010327" "博时消费创新混合C"     "2.29"  "80.27" "2.91" "0.0666" 10
Set-Q4Row 6 4 "009619" "博时女性消费主题混合A" "0.56"  "72.72" "9.83" "0.0550" 1
Set-Q4Row 7 5 "013837" "博时时代消费混合C"     "0.11"  "92.95" "6.26" "0.0069" 1
Set-Q4Row 8 6 "009620" "博时女性消费主题混合C" "0.03"  "72.72" "9.83" "0.0029" 1
Set-Q4Row 9 7 "011879" "博时新兴消费主题混合C" "0.00"  "87.10" "5.81" 0       5

# ---------------------------------------------------------------------------
# Update the "总计" sheet: insert the new 2022-Q4 row at the top of the data
# (row 2), pushing the existing rows down by one.
# ---------------------------------------------------------------------------

# Give the new row (row 5) the same per-row styling (the "A" marker column)
# as the existing last data row before we overwrite everything.
$wsTotal.Range("A4").Copy()
$wsTotal.Range("A5").PasteSpecial(-4122)

function Set-TotalRow($row, $idx, $label, $count, $mv) {
    $wsTotal.Cells.Item($row, 1).Value = $idx
    $wsTotal.Cells.Item($row, 2).Value = $label
    $wsTotal.Cells.Item($row, 3).Value = $count
    $wsTotal.Cells.Item($row, 4).Value = $mv
}

Set-TotalRow 2 0 "2022-Q4" 8 1.11
Set-TotalRow 3 1 "2022-Q3" 2 0.03
Set-TotalRow 4 2 "2021-Q2" 1 1.52
Set-TotalRow 5 3 "2021-Q1" 1 1.18

# Restore "总计" as the active sheet (as it was before the edit), since
# copying sheets around tends to leave a different sheet activated.
$wsTotal.Activate()
